# Auto-generated Excel COM-interop script
# Applies numeric value updates to cached cells across multiple worksheets
# as produced by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1632.1621
$ws.Range("J17").Value = 1632.1621
$ws.Range("L17").Value = 4896.4863
$ws.Range("N17").Value = -5232.4863
$ws.Range("H33").Value = 204.05
$ws.Range("I33").Value = 176.72223
$ws.Range("K33").Value = 176.72223
$ws.Range("M33").Value = 52.27777
$ws.Range("H116").Value = 40339.383
$ws.Range("I116").Value = 54182.355
$ws.Range("J116").Value = 12653.429
$ws.Range("K116").Value = 54182.355
$ws.Range("L116").Value = 12653.429
$ws.Range("M116").Value = -50740.355
$ws.Range("N116").Value = -19537.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4961.6523
$ws.Range("I61").Value = 2195.3333
$ws.Range("K61").Value = 2195.3333
$ws.Range("M61").Value = -1983.3333
$ws.Range("H97").Value = 791.44446
$ws.Range("I97").Value = 680.7273
$ws.Range("J97").Value = 1278.6
$ws.Range("K97").Value = 680.7273
$ws.Range("L97").Value = 1278.6
$ws.Range("M97").Value = -184.7273
$ws.Range("N97").Value = -2270.6
$ws.Range("H102").Value = 3237.375
$ws.Range("I102").Value = 2807
$ws.Range("K102").Value = 2807
$ws.Range("M102").Value = -1185
$ws.Range("H110").Value = 474.77777
$ws.Range("I110").Value = 478.16
$ws.Range("K110").Value = 478.16
$ws.Range("M110").Value = 1566.84
$ws.Range("H119").Value = 200558.2
$ws.Range("J119").Value = 200558.2
$ws.Range("L119").Value = 200558.2
$ws.Range("N119").Value = -210234.2
$ws.Range("H132").Value = 5757.3335
$ws.Range("I132").Value = 5757.3335
$ws.Range("K132").Value = 17272.0005
$ws.Range("M132").Value = -14742.0005
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
$ws.Range("H136").Value = 4961.6523
$ws.Range("I136").Value = 2195.3333
$ws.Range("K136").Value = 6585.999899999999
$ws.Range("M136").Value = -4035.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1909.3684
$ws.Range("I20").Value = 1895.0834
$ws.Range("K20").Value = 1895.0834
$ws.Range("M20").Value = -1648.0834
$ws.Range("H54").Value = 8142.75
$ws.Range("I54").Value = 7493
$ws.Range("K54").Value = 7493
$ws.Range("M54").Value = -7009
$ws.Range("H86").Value = 1356.3
$ws.Range("I86").Value = 1459.375
$ws.Range("J86").Value = 944
$ws.Range("K86").Value = 1459.375
$ws.Range("L86").Value = 944
$ws.Range("M86").Value = -336.375
$ws.Range("N86").Value = -3190
$ws.Range("H89").Value = 1356.3
$ws.Range("I89").Value = 1459.375
$ws.Range("J89").Value = 944
$ws.Range("K89").Value = 7296.875
$ws.Range("L89").Value = 4720
$ws.Range("M89").Value = -1680.875
$ws.Range("N89").Value = -15952
$ws.Range("H94").Value = 1077.1177
$ws.Range("I94").Value = 1121.1072
$ws.Range("J94").Value = 871.8333
$ws.Range("K94").Value = 1121.1072
$ws.Range("L94").Value = 871.8333
$ws.Range("M94").Value = -670.1071999999999
$ws.Range("N94").Value = -1773.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4651.613
$ws.Range("I31").Value = 2186.64
$ws.Range("J31").Value = 6317.1353
$ws.Range("K31").Value = 2186.64
$ws.Range("L31").Value = 6317.1353
$ws.Range("M31").Value = -1891.64
$ws.Range("N31").Value = -6907.1353
$ws.Range("H34").Value = 4651.613
$ws.Range("I34").Value = 2186.64
$ws.Range("J34").Value = 6317.1353
$ws.Range("K34").Value = 2186.64
$ws.Range("L34").Value = 6317.1353
$ws.Range("M34").Value = -1984.64
$ws.Range("N34").Value = -6721.1353
$ws.Range("H35").Value = 12902.625
$ws.Range("I35").Value = 14317.286
$ws.Range("K35").Value = 14317.286
$ws.Range("M35").Value = -14023.286
$ws.Range("H42").Value = 3000
$ws.Range("I42").Value = 3000
$ws.Range("K42").Value = 3000
$ws.Range("M42").Value = -2407
$ws.Range("H99").Value = 6627.125
$ws.Range("J99").Value = 6169.6665
$ws.Range("L99").Value = 6169.6665
$ws.Range("N99").Value = -9165.666499999999
$ws.Range("H105").Value = 1995
$ws.Range("I105").Value = 1990
$ws.Range("K105").Value = 1990
$ws.Range("M105").Value = -243
$ws.Range("H122").Value = 2460.1428
$ws.Range("I122").Value = 1877
$ws.Range("K122").Value = 5631
$ws.Range("M122").Value = -3181
$ws.Range("H126").Value = 6627.125
$ws.Range("J126").Value = 6169.6665
$ws.Range("L126").Value = 18508.9995
$ws.Range("N126").Value = -23448.9995
$ws.Range("H132").Value = 3004.8572
$ws.Range("I132").Value = 2402
$ws.Range("K132").Value = 7206
$ws.Range("M132").Value = -4676

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 2928711.8
$ws.Range("I7").Value = 1666830.5
$ws.Range("J7").Value = 10500000
$ws.Range("K7").Value = 5000491.5
$ws.Range("L7").Value = 31500000
$ws.Range("M7").Value = -5000379.5
$ws.Range("N7").Value = -31500224
$ws.Range("H40").Value = 1350.4117
$ws.Range("I40").Value = 159.66667
$ws.Range("J40").Value = 1999.909
$ws.Range("K40").Value = 638.66668
$ws.Range("L40").Value = 7999.636
$ws.Range("M40").Value = -569.66668
$ws.Range("N40").Value = -8137.636
$ws.Range("H80").Value = 3444.4443
$ws.Range("I80").Value = 3200.4
$ws.Range("J80").Value = 3749.5
$ws.Range("K80").Value = 9601.200000000001
$ws.Range("L80").Value = 11248.5
$ws.Range("M80").Value = -8665.200000000001
$ws.Range("N80").Value = -13120.5
$ws.Range("H83").Value = 3444.4443
$ws.Range("I83").Value = 3200.4
$ws.Range("J83").Value = 3749.5
$ws.Range("K83").Value = 28803.6
$ws.Range("L83").Value = 33745.5
$ws.Range("M83").Value = -24123.6
$ws.Range("N83").Value = -43105.5
$ws.Range("H92").Value = 760.5625
$ws.Range("J92").Value = 1999.5
$ws.Range("L92").Value = 5998.5
$ws.Range("N92").Value = -8494.5
$ws.Range("H109").Value = 3982.25
$ws.Range("I109").Value = 1720
$ws.Range("K109").Value = 5160
$ws.Range("M109").Value = -4120
$ws.Range("H133").Value = 4906.357
$ws.Range("I133").Value = 4854.3335
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 14563.0005
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -9503.000499999998
$ws.Range("N133").Value = -25120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 8250
$ws.Range("I6").Value = 7500
$ws.Range("J6").Value = 9000
$ws.Range("K6").Value = 7500
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = -7387
$ws.Range("N6").Value = -9226
$ws.Range("H16").Value = 8250
$ws.Range("I16").Value = 7500
$ws.Range("J16").Value = 9000
$ws.Range("K16").Value = 7500
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = -7250
$ws.Range("N16").Value = -9500
$ws.Range("H46").Value = 22714.857
$ws.Range("I46").Value = 18833.166
$ws.Range("J46").Value = 25626.125
$ws.Range("K46").Value = 18833.166
$ws.Range("L46").Value = 25626.125
$ws.Range("M46").Value = -18677.166
$ws.Range("N46").Value = -25938.125
$ws.Range("H102").Value = 1895.6207
$ws.Range("I102").Value = 1121.5454
$ws.Range("J102").Value = 4328.4287
$ws.Range("K102").Value = 1121.5454
$ws.Range("L102").Value = 4328.4287
$ws.Range("M102").Value = 500.4546
$ws.Range("N102").Value = -7572.4287
$ws.Range("H132").Value = 2876.889
$ws.Range("I132").Value = 2849.75
$ws.Range("K132").Value = 8549.25
$ws.Range("M132").Value = -6019.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 228175
$ws.Range("I62").Value = 451500
$ws.Range("J62").Value = 4850
$ws.Range("K62").Value = 451500
$ws.Range("L62").Value = 4850
$ws.Range("M62").Value = -450876
$ws.Range("N62").Value = -6098
$ws.Range("H65").Value = 228175
$ws.Range("I65").Value = 451500
$ws.Range("J65").Value = 4850
$ws.Range("K65").Value = 2257500
$ws.Range("L65").Value = 24250
$ws.Range("M65").Value = -2254380
$ws.Range("N65").Value = -6098
$ws.Range("H119").Value = 144558.4
$ws.Range("J119").Value = 144558.4
$ws.Range("L119").Value = 144558.4
$ws.Range("N119").Value = -154234.4
$ws.Range("H122").Value = 3864.6667
$ws.Range("J122").Value = 4478
$ws.Range("L122").Value = 13434
$ws.Range("N122").Value = -18334
$ws.Range("H132").Value = 4069.3704
$ws.Range("I132").Value = 2865.7
$ws.Range("K132").Value = 8597.099999999999
$ws.Range("M132").Value = -6067.099999999999
$ws.Range("H136").Value = 57568020
$ws.Range("I136").Value = 115129544
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 345388632
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -345386082
$ws.Range("N136").Value = -24600
